# The commit swaps the content of ppt/theme/theme1.xml (the "Integral" theme,
# used by the slide master / main presentation) and ppt/theme/theme2.xml (the
# "Office Theme", used only by the notes master). Both themes already share an
# identical font scheme and format scheme (fills/lines/effects) -- the only
# real content difference between the two files is the 12 color-scheme RGB
# values (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
#
# The exposed PowerPoint object model lets us reach the theme that backs the
# slide master (ppt/theme/theme1.xml) via Master.Theme.ThemeColorScheme, so we
# repoint its 12 colors at what used to be the "Office Theme" palette -- the
# net, visible effect of the swap for the deck's actual theme.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$theme = $m.Theme
$cs = $theme.ThemeColorScheme

# Target palette: the colors that used to live in ppt/theme/theme2.xml
# ("Office Theme"), expressed as VBA RGB() integers (r + g*256 + b*65536).
$officeColors = @{
    1  = 0        # dk1      000000
    2  = 16777215 # lt1      FFFFFF
    3  = 6968388  # dk2      44546A
    4  = 15132391 # lt2      E7E6E6
    5  = 13998939 # accent1  5B9BD5
    6  = 3243501  # accent2  ED7D31
    7  = 10855845 # accent3  A5A5A5
    8  = 49407    # accent4  FFC000
    9  = 12874308 # accent5  4472C4
    10 = 4697456  # accent6  70AD47
    11 = 12673797 # hlink    0563C1
    12 = 7491477  # folHlink 954F72
}

for ($i = 1; $i -le $cs.Count; $i++) {
    $cs.Colors($i).RGB = $officeColors[$i]
}
